$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.426.75"
$ws.Range("E2").Value = "  -3.12%  "

$ws.Range("D3").Value = "3.700.00"
$ws.Range("E3").Value = "  -3.60%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'596.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.87%  "

$ws.Range("D6").Value = "'165.48"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.73%  "

$ws.Range("D7").Value = "3.694.32"
$ws.Range("E7").Value = "  -3.76%  "

$ws.Range("E8").Value = "  -0.18%  "

$ws.Range("E9").Value = "  -0.24%  "

$ws.Range("E10").Value = "  -3.03%  "

$ws.Range("E11").Value = "  -3.87%  "

$ws.Range("E12").Value = "  -3.53%  "

$ws.Range("D13").Value = "'37.58"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.37%  "

$ws.Range("E14").Value = "  -4.95%  "

$ws.Range("D15").Value = "4.318.44"
$ws.Range("E15").Value = "  -3.52%  "

$ws.Range("D16").Value = "3.700.52"
$ws.Range("E16").Value = "  -3.60%  "

$ws.Range("D17").Value = "67.467.90"
$ws.Range("E17").Value = "  -3.26%  "

$ws.Range("D18").Value = "'17.58"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +6.29%  "

$ws.Range("D19").Value = "'7.16"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.64%  "

$ws.Range("E20").Value = "  -3.28%  "

$ws.Range("D21").Value = "'490.92"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.49%  "

$ws.Range("D22").Value = "'9.13"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.50%  "

$ws.Range("E23").Value = "  -1.40%  "

$ws.Range("D24").Value = "'85.96"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.25%  "

$ws.Range("E25").Value = "  -5.80%  "

$ws.Range("E26").Value = "  -2.15%  "

$ws.Range("E27").Value = "  -3.12%  "

$ws.Range("D28").Value = "'10.14"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.00%  "

$ws.Range("E29").Value = "  +0.02%  "

$ws.Range("E30").Value = "  -1.24%  "

$ws.Range("E31").Value = "  -6.30%  "

$ws.Range("D32").Value = "'31.55"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.93%  "

$ws.Range("E33").Value = "  -3.63%  "

$ws.Range("D34").Value = "3.838.61"
$ws.Range("E34").Value = "  -3.55%  "

$ws.Range("D35").Value = "'0.107"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.41%  "

$ws.Range("D36").Value = "3.641.24"
$ws.Range("E36").Value = "  -3.47%  "

$ws.Range("E37").Value = "  -0.11%  "

$ws.Range("D38").Value = "'0.994"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.73%  "

$ws.Range("E39").Value = "  -5.08%  "

$ws.Range("E40").Value = "  -6.70%  "

$ws.Range("D41").Value = "'0.321"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.68%  "

$ws.Range("D42").Value = "'433.40"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -10.22%  "

$ws.Range("D43").Value = "'48.60"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.22%  "

$ws.Range("E44").Value = "  -5.45%  "

$ws.Range("E45").Value = "  -6.30%  "

$ws.Range("D46").Value = "'8.37"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.38%  "

$ws.Range("E47").Value = "  -0.03%  "

$ws.Range("D48").Value = "'40.60"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.75%  "

$ws.Range("D49").Value = "'142.77"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.02%  "

$ws.Range("D50").Value = "2.754.57"
$ws.Range("E50").Value = "  -5.71%  "

$ws.Range("E51").Value = "  -3.34%  "
